$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 143.4675496666667
$ws.Range("H2").Value = 430.402649
$ws.Range("I2").Value = 0.2436371325027481
$ws.Range("J2").Value = 0.2436371325027482
$ws.Range("O2").Value = 0.9418062875790357
$ws.Range("P2").Value = 0.9418062875790357
$ws.Range("Q2").Value = 29.61863651610055
$ws.Range("R2").Value = 266.567728644905
$ws.Range("S2").Value = 0.2294589832788148
$ws.Range("T2").Value = 0.2294589832788149
$ws.Range("G3").Value = 143.4675496666667
$ws.Range("H3").Value = 430.402649
$ws.Range("I3").Value = 0.2436371325027481
$ws.Range("J3").Value = 0.2436371325027482
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 0.3333333333333333
$ws.Range("M3").Value = 0.01275633333333333
$ws.Range("N3").Value = 0.038269
$ws.Range("O3").Value = 0.05819371242096427
$ws.Range("P3").Value = 0.05819371242096427
$ws.Range("Q3").Value = 1.830119886064555
$ws.Range("R3").Value = 16.471078974581
$ws.Range("S3").Value = 0.01417814922393329
$ws.Range("T3").Value = 0.01417814922393329
$ws.Range("I4").Value = 0.2680684099784185
$ws.Range("J4").Value = 0.2680684099784185
$ws.Range("O4").Value = 0.9418062875790357
$ws.Range("P4").Value = 0.9418062875790357
$ws.Range("S4").Value = 0.2524685140189892
$ws.Range("T4").Value = 0.2524685140189892
$ws.Range("I5").Value = 0.2680684099784185
$ws.Range("J5").Value = 0.2680684099784185
$ws.Range("K5").Value = 1
$ws.Range("L5").Value = 0.3333333333333333
$ws.Range("M5").Value = 0.01275633333333333
$ws.Range("N5").Value = 0.038269
$ws.Range("O5").Value = 0.05819371242096427
$ws.Range("P5").Value = 0.05819371242096427
$ws.Range("Q5").Value = 2.013639394322111
$ws.Range("R5").Value = 18.122754548899
$ws.Range("S5").Value = 0.01559989595942923
$ws.Range("T5").Value = 0.01559989595942923
$ws.Range("G6").Value = 72.92931733333334
$ws.Range("H6").Value = 218.787952
$ws.Range("I6").Value = 0.1238488410219541
$ws.Range("J6").Value = 0.1238488410219541
$ws.Range("O6").Value = 0.9418062875790357
$ws.Range("P6").Value = 0.9418062875790357
$ws.Range("Q6").Value = 15.05613601460445
$ws.Range("R6").Value = 135.50522413144
$ws.Range("S6").Value = 0.1166416171838528
$ws.Range("T6").Value = 0.1166416171838528
$ws.Range("G7").Value = 72.92931733333334
$ws.Range("H7").Value = 218.787952
$ws.Range("I7").Value = 0.1238488410219541
$ws.Range("J7").Value = 0.1238488410219541
$ws.Range("K7").Value = 1
$ws.Range("L7").Value = 0.3333333333333333
$ws.Range("M7").Value = 0.01275633333333333
$ws.Range("N7").Value = 0.038269
$ws.Range("O7").Value = 0.05819371242096427
$ws.Range("P7").Value = 0.05819371242096427
$ws.Range("Q7").Value = 0.9303106816764446
$ws.Range("R7").Value = 8.372796135088
$ws.Range("S7").Value = 0.007207223838101319
$ws.Range("T7").Value = 0.007207223838101319
$ws.Range("G8").Value = 86.33190533333334
$ws.Range("H8").Value = 258.995716
$ws.Range("I8").Value = 0.1466091663779145
$ws.Range("J8").Value = 0.1466091663779145
$ws.Range("O8").Value = 0.9418062875790357
$ws.Range("P8").Value = 0.9418062875790357
$ws.Range("Q8").Value = 17.82307796955778
$ws.Range("R8").Value = 160.40770172602
$ws.Range("S8").Value = 0.1380774347114408
$ws.Range("T8").Value = 0.1380774347114408
$ws.Range("G9").Value = 86.33190533333334
$ws.Range("H9").Value = 258.995716
$ws.Range("I9").Value = 0.1466091663779145
$ws.Range("J9").Value = 0.1466091663779145
$ws.Range("K9").Value = 1
$ws.Range("L9").Value = 0.3333333333333333
$ws.Range("M9").Value = 0.01275633333333333
$ws.Range("N9").Value = 0.038269
$ws.Range("O9").Value = 0.05819371242096427
$ws.Range("P9").Value = 0.05819371242096427
$ws.Range("Q9").Value = 1.101278561733778
$ws.Range("R9").Value = 9.911507055604
$ws.Range("S9").Value = 0.008531731666473658
$ws.Range("T9").Value = 0.008531731666473658
$ws.Range("G10").Value = 34.55480166666666
$ws.Range("H10").Value = 103.664405
$ws.Range("I10").Value = 0.05868109416957502
$ws.Range("J10").Value = 0.05868109416957502
$ws.Range("O10").Value = 0.9418062875790357
$ws.Range("P10").Value = 0.9418062875790357
$ws.Range("Q10").Value = 7.133781212747222
$ws.Range("R10").Value = 64.204030914725
$ws.Range("S10").Value = 0.05526622345092325
$ws.Range("T10").Value = 0.05526622345092325
$ws.Range("G11").Value = 34.55480166666666
$ws.Range("H11").Value = 103.664405
$ws.Range("I11").Value = 0.05868109416957502
$ws.Range("J11").Value = 0.05868109416957502
$ws.Range("K11").Value = 1
$ws.Range("L11").Value = 0.3333333333333333
$ws.Range("M11").Value = 0.01275633333333333
$ws.Range("N11").Value = 0.038269
$ws.Range("O11").Value = 0.05819371242096427
$ws.Range("P11").Value = 0.05819371242096427
$ws.Range("Q11").Value = 0.4407925683272221
$ws.Range("R11").Value = 3.967133114944999
$ws.Range("S11").Value = 0.003414870718651772
$ws.Range("T11").Value = 0.003414870718651772
$ws.Range("G12").Value = 93.71982300000001
$ws.Range("H12").Value = 281.159469
$ws.Range("I12").Value = 0.1591553559493899
$ws.Range("J12").Value = 0.1591553559493899
$ws.Range("O12").Value = 0.9418062875790357
$ws.Range("P12").Value = 0.9418062875790357
$ws.Range("Q12").Value = 19.348301258645
$ws.Range("R12").Value = 174.134711327805
$ws.Range("S12").Value = 0.1498935149350149
$ws.Range("T12").Value = 0.1498935149350149
$ws.Range("G13").Value = 93.71982300000001
$ws.Range("H13").Value = 281.159469
$ws.Range("I13").Value = 0.1591553559493899
$ws.Range("J13").Value = 0.1591553559493899
$ws.Range("K13").Value = 1
$ws.Range("L13").Value = 0.3333333333333333
$ws.Range("M13").Value = 0.01275633333333333
$ws.Range("N13").Value = 0.038269
$ws.Range("O13").Value = 0.05819371242096427
$ws.Range("P13").Value = 0.05819371242096427
$ws.Range("Q13").Value = 1.195521302129
$ws.Range("R13").Value = 10.759691719161
$ws.Range("S13").Value = 0.009261841014374998
$ws.Range("T13").Value = 0.009261841014374998
